$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("types")
$ws.Range("A1").Value = "TEST"
